# Reorders the work-experience entries so that the sequence becomes:
#   Siege Analytics, Myers Research, PCCC, Feldman Group, Lake Research Partners
# (originally: Siege Analytics, PCCC, Myers Research, Lake Research Partners, Feldman Group)
#
# Strategy:
#   1. Insert a fresh copy of the "Myers Research" entry (heading + 4 body
#      paragraphs) immediately before the "Research Director - PCCC" heading.
#   2. Overwrite the text of the now-duplicated old "Myers Research" entry
#      (which sits right after PCCC) with the "Feldman Group" entry's content.
#   3. Delete the old "Feldman Group" entry, which now trails the
#      "Lake Research Partners" entry at the end of the section.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Step 1: insert the "Myers Research" block before "Research Director - PCCC"
# ---------------------------------------------------------------------

# Locate the PCCC heading paragraph explicitly.
$pcccParaIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text.StartsWith("Research Director - PCCC")) {
        $pcccParaIndex = $i
        break
    }
}

$insertRange = $d.Paragraphs($pcccParaIndex).Range
$insertRange.Collapse(1)

for ($i = 0; $i -lt 5; $i++) {
    $insertRange.InsertParagraphBefore()
}

$d.Paragraphs($pcccParaIndex).Range.Text = "Senior Analyst - Myers Research (Austin, TX) | 2012 - 2014"
$d.Paragraphs($pcccParaIndex).Style = "Heading3"

$d.Paragraphs($pcccParaIndex + 1).Range.Text = "Political Research & Analysis"
$d.Paragraphs($pcccParaIndex + 1).Style = "Normal"

$d.Paragraphs($pcccParaIndex + 2).Range.Text = "• Designed comprehensive survey instruments for specialized voting segments and niche markets"
$d.Paragraphs($pcccParaIndex + 2).Style = "Normal"

$d.Paragraphs($pcccParaIndex + 3).Range.Text = "• Developed sophisticated analytical products and reports that delivered actionable insights to clients"
$d.Paragraphs($pcccParaIndex + 3).Style = "Normal"

$d.Paragraphs($pcccParaIndex + 4).Range.Text = "• Co-developed a web application to manage all aspects of survey operations, from instrument design to data collection and analysis"
$d.Paragraphs($pcccParaIndex + 4).Style = "Normal"

# ---------------------------------------------------------------------
# Step 2: replace the old "Myers Research" block (now right after PCCC)
# with the "Feldman Group" entry's content
# ---------------------------------------------------------------------
$oldMyersIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text.StartsWith("Senior Analyst - Myers Research") -and $i -ne $pcccParaIndex) {
        $oldMyersIndex = $i
        break
    }
}

$d.Paragraphs($oldMyersIndex).Range.Text = "Field Director - The Feldman Group (Austin, TX) | 2011 - 2012"
$d.Paragraphs($oldMyersIndex + 1).Range.Text = "Political Campaign Management"
$d.Paragraphs($oldMyersIndex + 2).Range.Text = "• Managed all aspects of survey fielding for a multi-million dollar research firm, including scheduling, oversight, sampling, and quality control"
$d.Paragraphs($oldMyersIndex + 3).Range.Text = "• Developed and implemented data warehousing solutions for efficient storage and retrieval of research findings"
$d.Paragraphs($oldMyersIndex + 4).Range.Text = "• Created custom reports and data visualizations based on specific client requirements"

# ---------------------------------------------------------------------
# Step 3: delete the old "Feldman Group" block that used to trail
# "Lake Research Partners" at the end of the section
# ---------------------------------------------------------------------
$oldFeldmanIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text.StartsWith("Field Director - The Feldman Group") -and $i -ne $oldMyersIndex) {
        $oldFeldmanIndex = $i
        break
    }
}

$delStart = $d.Paragraphs($oldFeldmanIndex).Range.Start
$delEnd = $d.Paragraphs($oldFeldmanIndex + 4).Range.End
$delRange = $d.Range($delStart, $delEnd)
$delRange.Delete()
